# Update the "Förändrad" (Changed) date column (column C) for every data
# row on the sheet: bump the date serial value from 45179 to 45180
# (i.e. 2023-09-10 -> 2023-09-11) for all rows that currently hold 45179.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = $ws.UsedRange.Rows.Count }

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45179) {
        $cell.Value2 = 45180
    }
}
